$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.867074728012085
$ws.Range("B1").Value = 3.317430257797241
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.059309244155884
$ws.Range("E1").Value = 2.596883058547974
